# Added get_form_type in the import_utils
# This adds a new "Form Tag" column (column S) to the capital commitments
# worksheet, defaulting every existing data row to "Default".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column S
$ws.Range("S1").Value = "Form Tag"

# Populate the new column for every existing data row (rows 2-9) with "Default"
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 19).Value = "Default"
}

# Reflect the view state from the edit: scrolled so column M is the
# left-most visible column, with S3:S9 selected (active cell S3).
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 13
$ws.Range("S3:S9").Select()
